$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.479.17"
$ws.Range("E2").Value = "  -3.31%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.487.54"
$ws.Range("E3").Value = "  -0.81%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "552.41"
$ws.Range("E5").Value = "  -0.93%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.92"
$ws.Range("E6").Value = "  -6.47%  "
$ws.Range("E7").Value = "  +4.73%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -1.13%  "
$ws.Range("E10").Value = "  +1.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.63"
$ws.Range("E11").Value = "  -6.12%  "
$ws.Range("E12").Value = "  -1.88%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.20"
$ws.Range("E13").Value = "  -2.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.042.66"
$ws.Range("E14").Value = "  -1.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.488.75"
$ws.Range("E15").Value = "  -0.97%  "
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.40"
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.13"
$ws.Range("E18").Value = "  +2.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "65.468.09"
$ws.Range("E19").Value = "  -4.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.993"
$ws.Range("E20").Value = "  -1.65%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "412.13"
$ws.Range("E21").Value = "  +1.09%  "
$ws.Range("E22").Value = "  +1.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "85.62"
$ws.Range("E23").Value = "  +0.96%  "
$ws.Range("E24").Value = "  -3.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.65"
$ws.Range("E25").Value = "  +6.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.78"
$ws.Range("E26").Value = "  -8.39%  "
$ws.Range("E27").Value = "  -2.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.03"
$ws.Range("E28").Value = "  -1.99%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.02"
$ws.Range("E29").Value = "  +4.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.23"
$ws.Range("E30").Value = "  -1.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "612.73"
$ws.Range("E31").Value = "  -10.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.46"
$ws.Range("E32").Value = "  -6.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.65"
$ws.Range("E33").Value = "  -0.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.109"
$ws.Range("E34").Value = "  -1.46%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "59.42"
$ws.Range("E35").Value = "  -1.96%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.148"
$ws.Range("E36").Value = "  +11.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("B38").Value = "InjectiveProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "37.07"
$ws.Range("E38").Value = "  -5.46%  "
$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0787"
$ws.Range("E39").Value = "  -6.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.357.87"
$ws.Range("E40").Value = "  +9.55%  "
$ws.Range("E41").Value = "  -6.50%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.16%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.27"
$ws.Range("E43").Value = "  -4.35%  "
$ws.Range("E44").Value = "  -6.10%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.27"
$ws.Range("E45").Value = "  +3.18%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.52"
$ws.Range("E46").Value = "  -9.02%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0414"
$ws.Range("E47").Value = "  -2.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.74"
$ws.Range("E48").Value = "  -0.51%  "
$ws.Range("E49").Value = "  +1.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "137.89"
$ws.Range("E50").Value = "  -0.82%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.41"
$ws.Range("E51").Value = "  -10.15%  "
